$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 4th sheet (xl/worksheets/sheet4.xml, rId4).
# The diff shows a new blank column inserted before column N (Late/Outstanding
# data shifts from N/O/P to O/P/Q), which is exactly what Excel's
# "insert column" does.
$ws4 = $wb.Worksheets.Item("Repayment schedule")
$ws4.Columns.Item(14).Insert() | Out-Null

# The active sheet/selection also moved: "Edit Repayment Schedule1" (sheet3)
# was tabSelected with B4 selected; now "Repayment schedule" (sheet4) is the
# tabSelected sheet, with J20 selected.
$ws4.Activate()
$ws4.Range("J20").Select() | Out-Null

$wb.Save() | Out-Null
